$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1. Date change
Replace-Text "2024-10-17" "2024-10-18"

# 2a. First sentence of the "for assessing..." paragraph (single run edit)
Replace-Text "In terms of reliability and validity, research has shown that CJ can generate precise and consistent scores that accurately represent the traits being measured. Notable contributions in this research area include studies by Pollit" "Research on reliability and validity has shown that CJ can generate precise and consistent scores that accurately represent the traits being measured"

# 2b. Merge citation list into a single run (Pollitt et al.)
Replace-Text "(2012a, 2012b), Whitehouse (2012), van Daal et al. (2016), Lesterhuis (2018), Bramley and Vitello (2019), Verhavert et al. (2019), Crompvoets, Béguin, and Sijtsma (2022), and Bouwer et al. (2023)" "(Pollitt 2012a, 2012b; Whitehouse 2012; van Daal et al. 2016; Lesterhuis 2018; Bramley and Vitello 2019; Verhavert et al. 2019; Crompvoets, Béguin, and Sijtsma 2022; Bouwer et al. 2023)"

# 2c. Second sentence (practical applicability intro) - single run edit
Replace-Text ". Regarding practical applicability, several studies have highlighted the method’s versatility in both educational and non-educational contexts, presenting it as an efficient and effective alternative for measurement and evaluation. Key examples in this research area include the works of" ". Regarding practical applicability, several studies have highlighted CJ’s versatility across both educational and non-educational contexts, presenting it as an efficient and effective alternative for measurement and evaluation"

# 2d. Merge second citation list into a single run (Jones et al.)
Replace-Text "Jones (2015), Bartholomew et al. (2018), Jones et al. (2019), Marshall et al. (2020), Bartholomew and Williams (2020), and Boonen, Kloots, and Gillis (2020)" "(Jones 2015; Bartholomew et al. 2018; Jones et al. 2019; Marshall et al. 2020; Bartholomew and Williams 2020; Boonen, Kloots, and Gillis 2020)"

# 3. "concerning CJ..." paragraph
Replace-Text "Despite the growing number of CJ studies, the unsystematic and fragmented research approaches employed in the literature have overlooked several critical issues concerning CJ. These issues can be categorized into concerns about the method’s structural model, measurement model, and experimental design." "Despite the growing number of CJ studies, the unsystematic and fragmented research approaches employed in the literature have overlooked several critical issues concerning the method. These issues fall into three main categories: concerns about the structural model, the measurement model, and the experimental design of CJ. In the following sections, each issue will be discussed in detail, followed by the introduction of an approach that addresses all three concerns simultaneously."

# 4. "A common practice..." sentence
Replace-Text "A common practice in CJ literature involves performing data analysis and hypothesis testing on scores previously estimated using the Bradley-Terry-Luce (BTL) model" "In CJ literature, it is common to perform data analysis and hypothesis testing on scores estimated using the Bradley-Terry-Luce (BTL) model"

# 5. "These studies use the scores..."
Replace-Text ". These studies use the scores or their transformations to identify" ". Several studies use the scores generated by the BTL model or their transformations to identify"

# 6. "or test various hypothesis"
Replace-Text ", or test various hypothesis" ", or test various hypotheses about the underlying trait being measured"

# 7. "since CJ scores are parameter estimates"
Replace-Text ". However, since CJ scores are parameter estimates with inherent uncertainty" ". However, since these scores are parameter estimates with inherent uncertainty"

# 8. Final sentence addition after (McElreath 2020)
Replace-Text "(McElreath 2020)." "(McElreath 2020). To address this issue properly, the approach should follow a strategy similar to that used in Structural Equation Modeling (SEM), where data analysis and hypothesis testing occur at the structural model level, while the BTL model functions as the measurement model."
